$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "LOB1207 -  Poluição Ambiental I  (Requisito fraco)" requirement row
# (row 25). This shifts the following row ("LOQ4233 ...", previously row 26) up
# to become the new last row, and the worksheet shrinks from 26 to 25 rows.
$ws.Rows.Item(25).Delete()
